# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   across Overview (E/F) and the per-locale sheets (Status column C).
# - The per-locale "Latest Handback DateTime" (col K) is refreshed to the
#   handback-generation timestamp.
# - The stale "handback file is not the latest" Error Detail (col P) is
#   cleared now that the report is back in sync.
# - Columns that now hold the longer status text / shorter (empty) error
#   text are resized to fit their new content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both data rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), Latest Handback DateTime (K), Error Detail (P) ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-11-03 20:27:28"
$wsZhCn.Range("K3").Value = "2016-11-03 20:27:28"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: Status column (C), Latest Handback DateTime (K), Error Detail (P) ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-11-03 20:27:47"
$wsDeDe.Range("K3").Value = "2016-11-03 20:27:47"
$wsDeDe.Range("P2").Value = ""

# --- Resize columns to fit the new content (status text widened, error column shrunk) ---
$wsOverview.Range("E1").ColumnWidth = 29.14437166849777
$wsOverview.Range("F1").ColumnWidth = 29.14437166849777

$wsZhCn.Range("C1").ColumnWidth = 29.14437166849777
$wsZhCn.Range("P1").ColumnWidth = 12.913719813028965

$wsDeDe.Range("C1").ColumnWidth = 29.14437166849777
$wsDeDe.Range("P1").ColumnWidth = 12.913719813028965
